# Generate Report for Handback
# The localization CI run completed: the de-de / zh-cn handback packages are
# now in sync with en-US, so the per-language "Status" flips from
# "Ready for handoff" to "Handed back: in sync with en-US", the stale
# "handback file is not the latest" error clears, and the Latest Handback
# DateTime stamps advance to the moment the report was (re)generated.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: per-language status columns ---
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$ws2.Range("C2").Value = $newStatus
$ws2.Range("K2").Value = "2016-08-15 18:44:38"
$ws2.Range("P2").ClearContents()

# --- de-de sheet ---
$ws3.Range("C2").Value = $newStatus
$ws3.Range("K2").Value = "2016-08-15 18:44:45"
$ws3.Range("P2").ClearContents()

# --- Column widths: the longer status text no longer fits the old
# column width, so the Status columns widen and the now-empty Error
# Detail columns shrink back down (mirrors Excel's own autofit pass). ---
$ws1.Columns.Item(5).ColumnWidth = 29.16666666666667
$ws1.Columns.Item(6).ColumnWidth = 29.16666666666667

$ws2.Columns.Item(3).ColumnWidth = 29.16666666666667
$ws2.Columns.Item(16).ColumnWidth = 12.83333333333333

$ws3.Columns.Item(3).ColumnWidth = 29.16666666666667
$ws3.Columns.Item(16).ColumnWidth = 12.83333333333333
